# Auto-generated script applying Universalis price refresh updates
# to the Mandragora_Profits workbook (columns H-N: price/profit data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 220.25
$ws.Range("I9").Value = 223.4
$ws.Range("J9").Value = 215
$ws.Range("K9").Value = 223.4
$ws.Range("L9").Value = 215
$ws.Range("M9").Value = -54.40000000000001
$ws.Range("N9").Value = -553
$ws.Range("H19").Value = 488.30304
$ws.Range("I19").Value = 391.53845
$ws.Range("J19").Value = 551.2
$ws.Range("K19").Value = 391.53845
$ws.Range("L19").Value = 551.2
$ws.Range("M19").Value = -216.53845
$ws.Range("N19").Value = -901.2
$ws.Range("H137").Value = 2022.2
$ws.Range("I137").Value = 2179.35
$ws.Range("J137").Value = 1812.6666
$ws.Range("K137").Value = 6538.049999999999
$ws.Range("L137").Value = 5437.9998
$ws.Range("M137").Value = -3988.049999999999
$ws.Range("N137").Value = -10537.9998
$ws.Range("H138").Value = 2371.6223
$ws.Range("I138").Value = 1172.7667
$ws.Range("J138").Value = 4769.3335
$ws.Range("K138").Value = 3518.300099999999
$ws.Range("L138").Value = 14308.0005
$ws.Range("M138").Value = 1621.699900000001
$ws.Range("N138").Value = -24588.0005
$ws.Range("H141").Value = 6770.689
$ws.Range("I141").Value = 1177.6216
$ws.Range("J141").Value = 32638.625
$ws.Range("K141").Value = 3532.8648
$ws.Range("L141").Value = 97915.875
$ws.Range("M141").Value = 1647.1352
$ws.Range("N141").Value = -108275.875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2769.889
$ws.Range("I2").Value = 2431.75
$ws.Range("J2").Value = 5475
$ws.Range("K2").Value = 2431.75
$ws.Range("L2").Value = 5475
$ws.Range("M2").Value = -2318.75
$ws.Range("N2").Value = -5701
$ws.Range("H32").Value = 4957.193
$ws.Range("I32").Value = 4414.5396
$ws.Range("J32").Value = 10848.857
$ws.Range("K32").Value = 4414.5396
$ws.Range("L32").Value = 10848.857
$ws.Range("M32").Value = -4127.5396
$ws.Range("N32").Value = -11422.857
$ws.Range("H37").Value = 19000
$ws.Range("I37").Value = 19000
$ws.Range("K37").Value = 19000
$ws.Range("M37").Value = -18727
$ws.Range("H44").Value = 47000
$ws.Range("J44").Value = 47000
$ws.Range("L44").Value = 47000
$ws.Range("N44").Value = -47976
$ws.Range("H55").Value = 41999
$ws.Range("J55").Value = 41999
$ws.Range("L55").Value = 41999
$ws.Range("N55").Value = -42629
$ws.Range("H61").Value = 2705.0244
$ws.Range("I61").Value = 2064.6667
$ws.Range("J61").Value = 3074.4614
$ws.Range("K61").Value = 2064.6667
$ws.Range("L61").Value = 3074.4614
$ws.Range("M61").Value = -1852.6667
$ws.Range("N61").Value = -3498.4614
$ws.Range("H63").Value = 7499.1665
$ws.Range("J63").Value = 13500
$ws.Range("L63").Value = 13500
$ws.Range("N63").Value = -14872
$ws.Range("H66").Value = 7499.1665
$ws.Range("J66").Value = 13500
$ws.Range("L66").Value = 67500
$ws.Range("N66").Value = -74364
$ws.Range("H74").Value = 1427.1915
$ws.Range("I74").Value = 1158.3334
$ws.Range("J74").Value = 1790.15
$ws.Range("K74").Value = 1158.3334
$ws.Range("L74").Value = 1790.15
$ws.Range("M74").Value = -284.3334
$ws.Range("N74").Value = -3538.15
$ws.Range("H77").Value = 1427.1915
$ws.Range("I77").Value = 1158.3334
$ws.Range("J77").Value = 1790.15
$ws.Range("K77").Value = 5791.666999999999
$ws.Range("L77").Value = 8950.75
$ws.Range("M77").Value = -1423.666999999999
$ws.Range("N77").Value = -17686.75
$ws.Range("H80").Value = 44000
$ws.Range("J80").Value = 44000
$ws.Range("L80").Value = 44000
$ws.Range("N80").Value = -45996
$ws.Range("H83").Value = 44000
$ws.Range("J83").Value = 44000
$ws.Range("L83").Value = 132000
$ws.Range("N83").Value = -141984
$ws.Range("H116").Value = 2769.889
$ws.Range("I116").Value = 2431.75
$ws.Range("J116").Value = 5475
$ws.Range("K116").Value = 2431.75
$ws.Range("L116").Value = 5475
$ws.Range("M116").Value = -137.75
$ws.Range("N116").Value = -10063
$ws.Range("H122").Value = 1426.5769
$ws.Range("I122").Value = 1218.619
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 3655.857
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -1205.857
$ws.Range("N122").Value = -11800
$ws.Range("H136").Value = 2705.0244
$ws.Range("I136").Value = 2064.6667
$ws.Range("J136").Value = 3074.4614
$ws.Range("K136").Value = 6194.000100000001
$ws.Range("L136").Value = 9223.3842
$ws.Range("M136").Value = -3644.000100000001
$ws.Range("N136").Value = -14323.3842
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2769.889
$ws.Range("I3").Value = 2431.75
$ws.Range("J3").Value = 5475
$ws.Range("K3").Value = 2431.75
$ws.Range("L3").Value = 5475
$ws.Range("M3").Value = -2317.75
$ws.Range("N3").Value = -5703
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H35").Value = 53815.832
$ws.Range("J35").Value = 55162.727
$ws.Range("L35").Value = 55162.727
$ws.Range("N35").Value = -55782.727
$ws.Range("H134").Value = 4966.531
$ws.Range("I134").Value = 2602.5217
$ws.Range("J134").Value = 7057.769
$ws.Range("K134").Value = 7807.5651
$ws.Range("L134").Value = 21173.307
$ws.Range("M134").Value = -5272.5651
$ws.Range("N134").Value = -26243.307
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7094417
$ws.Range("I31").Value = 1669.9546
$ws.Range("J31").Value = 13336035
$ws.Range("K31").Value = 1669.9546
$ws.Range("L31").Value = 13336035
$ws.Range("M31").Value = -1374.9546
$ws.Range("N31").Value = -13336625
$ws.Range("H34").Value = 7094417
$ws.Range("I34").Value = 1669.9546
$ws.Range("J34").Value = 13336035
$ws.Range("K34").Value = 1669.9546
$ws.Range("L34").Value = 13336035
$ws.Range("M34").Value = -1467.9546
$ws.Range("N34").Value = -13336439
$ws.Range("H132").Value = 2495.6978
$ws.Range("I132").Value = 1790.9565
$ws.Range("J132").Value = 3306.15
$ws.Range("K132").Value = 5372.8695
$ws.Range("L132").Value = 9918.450000000001
$ws.Range("M132").Value = -2842.8695
$ws.Range("N132").Value = -14978.45
$ws.Range("H134").Value = 1039.6123
$ws.Range("I134").Value = 449.5
$ws.Range("J134").Value = 1826.4286
$ws.Range("K134").Value = 1348.5
$ws.Range("L134").Value = 5479.2858
$ws.Range("M134").Value = 1186.5
$ws.Range("N134").Value = -10549.2858
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2780.5642
$ws.Range("J122").Value = 3070.9697
$ws.Range("L122").Value = 27638.7273
$ws.Range("N122").Value = -32538.7273
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H101").Value = 30109.5
$ws.Range("J101").Value = 30109.5
$ws.Range("L101").Value = 30109.5
$ws.Range("N101").Value = -36599.5
$ws.Range("H109").Value = 20000
$ws.Range("J109").Value = 20000
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22080
$ws.Range("H132").Value = 1897066.1
$ws.Range("I132").Value = 3790233
$ws.Range("J132").Value = 3899.2727
$ws.Range("K132").Value = 11370699
$ws.Range("L132").Value = 11697.8181
$ws.Range("M132").Value = -11368169
$ws.Range("N132").Value = -16757.8181
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1714.3158
$ws.Range("I7").Value = 1369.0714
$ws.Range("J7").Value = 2681
$ws.Range("K7").Value = 1369.0714
$ws.Range("L7").Value = 2681
$ws.Range("M7").Value = -1257.0714
$ws.Range("N7").Value = -2905
$ws.Range("H40").Value = 3905.12
$ws.Range("I40").Value = 2949.8948
$ws.Range("J40").Value = 6930
$ws.Range("K40").Value = 2949.8948
$ws.Range("L40").Value = 6930
$ws.Range("M40").Value = -2813.8948
$ws.Range("N40").Value = -7202
$ws.Range("H126").Value = 1714.3158
$ws.Range("I126").Value = 1369.0714
$ws.Range("J126").Value = 2681
$ws.Range("K126").Value = 4107.2142
$ws.Range("L126").Value = 8043
$ws.Range("M126").Value = -1637.2142
$ws.Range("N126").Value = -12983
$ws.Range("H132").Value = 31253832
$ws.Range("I132").Value = 41670360
$ws.Range("J132").Value = 4249.125
$ws.Range("K132").Value = 125011080
$ws.Range("L132").Value = 12747.375
$ws.Range("M132").Value = -125008550
$ws.Range("N132").Value = -17807.375
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4632.657
$ws.Range("I122").Value = 3636.7273
$ws.Range("K122").Value = 10910.1819
$ws.Range("M122").Value = -8460.1819
$ws.Range("H124").Value = 31685.8
$ws.Range("J124").Value = 31685.8
$ws.Range("L124").Value = 31685.8
$ws.Range("N124").Value = -41505.8
$ws.Range("H125").Value = 29666.666
$ws.Range("J125").Value = 29666.666
$ws.Range("L125").Value = 29666.666
$ws.Range("N125").Value = -39506.666
$ws.Range("H132").Value = 1888.836
$ws.Range("I132").Value = 1456.3256
$ws.Range("J132").Value = 2922.0557
$ws.Range("K132").Value = 4368.976799999999
$ws.Range("L132").Value = 8766.167099999999
$ws.Range("M132").Value = -1838.976799999999
$ws.Range("N132").Value = -13826.1671
